$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 55
$ws.Range("H2").Value = 55
$ws.Range("E5").Value = 86
$ws.Range("F5").Value = 57
$ws.Range("H5").Value = 57
$ws.Range("F6").Value = 18
$ws.Range("H6").Value = 18
$ws.Range("E10").Value = 272
$ws.Range("F10").Value = 141
$ws.Range("H10").Value = 141
$ws.Range("E11").Value = 198
$ws.Range("F11").Value = 123
$ws.Range("H11").Value = 123
$ws.Range("E12").Value = 292
$ws.Range("F12").Value = 178
$ws.Range("H12").Value = 178
$ws.Range("F13").Value = 46
$ws.Range("H13").Value = 46
$ws.Range("F14").Value = 41
$ws.Range("H14").Value = 41
$ws.Range("F15").Value = 43
$ws.Range("H15").Value = 43
$ws.Range("E16").Value = 117
$ws.Range("F17").Value = 30
$ws.Range("H17").Value = 30
$ws.Range("E22").Value = 109
$ws.Range("F22").Value = 66
$ws.Range("H22").Value = 66
$ws.Range("E23").Value = 115
$ws.Range("F24").Value = 69
$ws.Range("H24").Value = 69
$ws.Range("E25").Value = 141
$ws.Range("F25").Value = 72
$ws.Range("H25").Value = 72
$ws.Range("E27").Value = 185
$ws.Range("F27").Value = 105
$ws.Range("H27").Value = 105
$ws.Range("F28").Value = 42
$ws.Range("H28").Value = 42
$ws.Range("F29").Value = 74
$ws.Range("H29").Value = 74
$ws.Range("F31").Value = 26
$ws.Range("H31").Value = 26
$ws.Range("F32").Value = 68
$ws.Range("H32").Value = 68
$ws.Range("E33").Value = 176
$ws.Range("F33").Value = 97
$ws.Range("H33").Value = 97
$ws.Range("E34").Value = 132
$ws.Range("F35").Value = 61
$ws.Range("H35").Value = 61
$ws.Range("E36").Value = 36
$ws.Range("E37").Value = 95
$ws.Range("F37").Value = 55
$ws.Range("H37").Value = 55
$ws.Range("F38").Value = 46
$ws.Range("H38").Value = 46
$ws.Range("E39").Value = 123
$ws.Range("F39").Value = 66
$ws.Range("H39").Value = 66
$ws.Range("E40").Value = 170
$ws.Range("F40").Value = 84
$ws.Range("H40").Value = 84
$ws.Range("E41").Value = 232
$ws.Range("F41").Value = 113
$ws.Range("H41").Value = 113
$ws.Range("E42").Value = 221
$ws.Range("E43").Value = 68
$ws.Range("E44").Value = 178
$ws.Range("F44").Value = 104
$ws.Range("H44").Value = 104
$ws.Range("E46").Value = 171
$ws.Range("F46").Value = 109
$ws.Range("H46").Value = 109
$ws.Range("E47").Value = 276
$ws.Range("F47").Value = 146
$ws.Range("H47").Value = 146
$ws.Range("E48").Value = 129
$ws.Range("F48").Value = 56
$ws.Range("H48").Value = 56
$ws.Range("E49").Value = 154
$ws.Range("F49").Value = 76
$ws.Range("H49").Value = 76
$ws.Range("E50").Value = 131
$ws.Range("F50").Value = 57
$ws.Range("H50").Value = 57
$ws.Range("E51").Value = 127
$ws.Range("F51").Value = 61
$ws.Range("H51").Value = 61
$ws.Range("E52").Value = 16
